$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice date from 2021-04-05 to 2021-04-06
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

# Update D/E column values for rows 2-13
$ws.Range("D2").Value = 0.03022767169916129
$ws.Range("E2").Value = 0.002853067047075664

$ws.Range("D3").Value = 0.02385899793347262
$ws.Range("E3").Value = 0.004294917680744748

$ws.Range("D4").Value = 0.05185323089158646
$ws.Range("E4").Value = 0.003573981415296634

$ws.Range("D5").Value = 0.1355939204050834
$ws.Range("E5").Value = -0.002870167145027747

$ws.Range("D6").Value = 0.03007989821556251
$ws.Range("E6").Value = -0.0030120481927709

$ws.Range("D7").Value = 0.120740989542119
$ws.Range("E7").Value = 0.00458657153777553

$ws.Range("D8").Value = 0.1014116122583677
$ws.Range("E8").Value = -0.002075863370447184

$ws.Range("D9").Value = 0.02790956601956012
$ws.Range("E9").Value = 0.001106684373616718

$ws.Range("D10").Value = 0.12088355461605
$ws.Range("E10").Value = -0.001202404809619195

$ws.Range("D11").Value = 0.2533140323398635
$ws.Range("E11").Value = -0.002362133187971183

$ws.Range("D12").Value = 0.1041265260791734
$ws.Range("E12").Value = 0.004299394176275229

$ws.Range("E13").Value = -0.00002761668382011973

# Restore sheet protection to match the original protected state
$ws.Protect()
